$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory labels to pluralized/"(s)" form for a handful of rows
$ws.Range("H12").Value = "line graph(s)"
$ws.Range("H15").Value = "line graph(s)"
$ws.Range("H16").Value = "line graph(s)"
$ws.Range("H17").Value = "line graph(s)"
$ws.Range("H18").Value = "line graph(s)"
$ws.Range("H23").Value = "drawing(s)"

# Remove the now-unused "is_viewed" column (column I) entirely,
# shrinking the sheet's used range back down to column H.
$ws.Columns.Item(9).Delete()
